$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the sheet. It belongs right above
# the existing row for "Provincia de Linares" / 44586, so insert a fresh
# row at position 95 (this pushes the former rows 95-185 down to 96-186).
$ws.Rows.Item(95).Insert()

# Fill in the newly inserted row 95 with the new record's data.
$ws.Cells.Item(95, 1).Value = 3
$ws.Cells.Item(95, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44587
$ws.Cells.Item(95, 5).Value = 5
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100101
$ws.Cells.Item(95, 8).Value = "Berries"
$ws.Cells.Item(95, 9).Value = 100101001
$ws.Cells.Item(95, 10).Value = "Arándano (blue)"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 200
$ws.Cells.Item(95, 14).Value = 4000
$ws.Cells.Item(95, 15).Value = 4500
$ws.Cells.Item(95, 16).Value = 4300
$ws.Cells.Item(95, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Linares"
$ws.Cells.Item(95, 19).Value = 2150
$ws.Cells.Item(95, 20).Value = 2
